# Update test data for importing
#
# Rows 2-11, columns W:Z used to hold volatile RANDBETWEEN() formulas.
# Replace them with frozen (static) values, formatted with a plain
# "0.00" (two decimal place) number format, matching a fresh "paste
# values" of a random snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (W, X, Y, Z) static values captured from the formulas.
$data = @{
    2  = @(961, 508, 939, 523)
    3  = @(380, 942, 697, 540)
    4  = @(999, 288, 107, 587)
    5  = @(394, 948, 342, 734)
    6  = @(303, 885, 109, 338)
    7  = @(572, 448, 777, 122)
    8  = @(257, 768, 723, 180)
    9  = @(171, 518, 273, 442)
    10 = @(371, 840, 881, 266)
    11 = @(614.35, 103, 255, 866.25)
}

# Columns W, X, Y, Z are columns 23-26.
$cols = @(23, 24, 25, 26)

foreach ($row in 2..11) {
    $values = $data[$row]
    for ($i = 0; $i -lt 4; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value = $values[$i]
    }
}

# Apply a plain 2-decimal-place number format (numFmtId 2) to the whole
# block in one go so a single new style is reused, like Excel would.
$ws.Range("W2:Z11").NumberFormat = "0.00"

# Scroll the view over and move the selection, matching the saved
# window state (top-left corner near column L, active cell Z12).
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Z12").Select()
